$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" placeholders: 18-11-2024 -> 19-11-2024
#    These live on the slide master and on every slide layout.
# ---------------------------------------------------------------------
$oldDate = "18-11-2024"
$newDate = "19-11-2024"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDate = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDate = $true
                }
            } catch {
                $isDate = $false
            }
            if ($isDate) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout
}

# ---------------------------------------------------------------------
# 2) Slide 1, "TextBox 9": widen the box and extend the sentence with
#    the new fabric-roll / tube / box label wording.
#    "... QR ren ..."  ->  "... QR Ren, vải hem, vải ống dạng thùng ..."
# ---------------------------------------------------------------------
$s = $ppt.ActivePresentation.Slides.Item(1)
$tb = $s.Shapes.Item(6)

# Widen the textbox to fit the longer caption.
$tb.Width = 6167073 / 914400 * 72

$tr = $tb.TextFrame.TextRange

# Grow " QR ren " into the full new phrase, then re-split it into the
# same word-by-word run boundaries the source deck ends up with.
$c = $tr.Characters(12, 8)
$c.Text = " QR Ren, vải hem, vải ống dạng thùng "

$c = $tr.Characters(12, 9)
$c.Text = " QR Ren, "

$c = $tr.Characters(21, 3)
$c.Text = "vải"

$c = $tr.Characters(24, 6)
$c.Text = " hem, "

$c = $tr.Characters(30, 3)
$c.Text = "vải"

$c = $tr.Characters(33, 1)
$c.Text = " "

$c = $tr.Characters(34, 3)
$c.Text = "ống"

$c = $tr.Characters(37, 1)
$c.Text = " "

$c = $tr.Characters(38, 4)
$c.Text = "dạng"

$c = $tr.Characters(42, 1)
$c.Text = " "

$c = $tr.Characters(43, 5)
$c.Text = "thùng"

$c = $tr.Characters(48, 1)
$c.Text = " "
